$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row (3 through 12), copy the value already present in column E
# out across columns F through K (same shared-string value repeated).
for ($r = 3; $r -le 12; $r++) {
    $val = $ws.Cells.Item($r, 5).Value2
    for ($c = 6; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).Value = $val
    }
}

# Update the active selection to match the new extent used (J14).
$ws.Range("J14").Select()
